$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1842
$ws.Range("I98").Value = 2258.8
$ws.Range("J98").Value = 800
$ws.Range("K98").Value = 2258.8
$ws.Range("L98").Value = 800
$ws.Range("M98").Value = -760.8000000000002
$ws.Range("N98").Value = -3796
$ws.Range("H116").Value = 2497.2727
$ws.Range("I116").Value = 1920
$ws.Range("J116").Value = 2978.3333
$ws.Range("K116").Value = 1920
$ws.Range("L116").Value = 2978.3333
$ws.Range("M116").Value = 1522
$ws.Range("N116").Value = -9862.3333
$ws.Range("H122").Value = 1842
$ws.Range("I122").Value = 2258.8
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 6776.400000000001
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = -4326.400000000001
$ws.Range("N122").Value = -7300

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 54004.332
$ws.Range("J23").Value = 46003.5
$ws.Range("L23").Value = 46003.5
$ws.Range("N23").Value = -46521.5
$ws.Range("H37").Value = 12700
$ws.Range("J37").Value = 14475
$ws.Range("L37").Value = 14475
$ws.Range("N37").Value = -15021
$ws.Range("H44").Value = 15857.143
$ws.Range("J44").Value = 15857.143
$ws.Range("L44").Value = 15857.143
$ws.Range("N44").Value = -16833.143
$ws.Range("H55").Value = 10110
$ws.Range("J55").Value = 10268.571
$ws.Range("L55").Value = 10268.571
$ws.Range("N55").Value = -10898.571
$ws.Range("H122").Value = 5307
$ws.Range("I122").Value = 6049.25
$ws.Range("J122").Value = 4119.4
$ws.Range("K122").Value = 18147.75
$ws.Range("L122").Value = 12358.2
$ws.Range("M122").Value = -15697.75
$ws.Range("N122").Value = -17258.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1204.3077
$ws.Range("I99").Value = 809.3889
$ws.Range("K99").Value = 809.3889
$ws.Range("M99").Value = 688.6111
$ws.Range("H107").Value = 15154135
$ws.Range("I107").Value = 37038496
$ws.Range("J107").Value = 3423.923
$ws.Range("K107").Value = 37038496
$ws.Range("L107").Value = 3423.923
$ws.Range("M107").Value = -37036576
$ws.Range("N107").Value = -7263.923

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 526
$ws.Range("I3").Value = 52
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 52
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = 61
$ws.Range("N3").Value = -1226
$ws.Range("H7").Value = 108.666664
$ws.Range("I7").Value = 39.5
$ws.Range("J7").Value = 247
$ws.Range("K7").Value = 39.5
$ws.Range("L7").Value = 247
$ws.Range("M7").Value = 73.5
$ws.Range("N7").Value = -473
$ws.Range("H31").Value = 15058.16
$ws.Range("I31").Value = 17785.355
$ws.Range("J31").Value = 5001.625
$ws.Range("K31").Value = 17785.355
$ws.Range("L31").Value = 5001.625
$ws.Range("M31").Value = -17490.355
$ws.Range("N31").Value = -5591.625
$ws.Range("H34").Value = 15058.16
$ws.Range("I34").Value = 17785.355
$ws.Range("J34").Value = 5001.625
$ws.Range("K34").Value = 17785.355
$ws.Range("L34").Value = 5001.625
$ws.Range("M34").Value = -17583.355
$ws.Range("N34").Value = -5405.625
$ws.Range("H63").Value = 41810.25
$ws.Range("J63").Value = 41810.25
$ws.Range("L63").Value = 41810.25
$ws.Range("N63").Value = -43182.25
$ws.Range("H66").Value = 41810.25
$ws.Range("J66").Value = 41810.25
$ws.Range("L66").Value = 125430.75
$ws.Range("N66").Value = -132294.75
$ws.Range("H80").Value = 10327.75
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 10327.75
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 10327.75
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -12573.75
$ws.Range("H83").Value = 10327.75
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 10327.75
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 30983.25
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -42215.25
$ws.Range("H133").Value = 45050
$ws.Range("J133").Value = 45050
$ws.Range("L133").Value = 45050
$ws.Range("N133").Value = -50110
$ws.Range("H135").Value = 48750
$ws.Range("J135").Value = 48750
$ws.Range("L135").Value = 48750
$ws.Range("N135").Value = -58890
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280
$ws.Range("H140").Value = 49800
$ws.Range("J140").Value = 49800
$ws.Range("L140").Value = 49800
$ws.Range("N140").Value = -60160

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 72.5
$ws.Range("I2").Value = 45
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 270
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = -157
$ws.Range("N2").Value = -826
$ws.Range("H37").Value = 1022506.1
$ws.Range("J37").Value = 1022506.1
$ws.Range("L37").Value = 3067518.3
$ws.Range("N37").Value = -3067742.3
$ws.Range("H58").Value = 1091.6666
$ws.Range("J58").Value = 1091.6666
$ws.Range("L58").Value = 3274.9998
$ws.Range("N58").Value = -3530.9998
$ws.Range("H103").Value = 1291.3334
$ws.Range("I103").Value = 249
$ws.Range("J103").Value = 1812.5
$ws.Range("K103").Value = 747
$ws.Range("L103").Value = 5437.5
$ws.Range("M103").Value = 132
$ws.Range("N103").Value = -7195.5
$ws.Range("H113").Value = 645.7368
$ws.Range("I113").Value = 613.4167
$ws.Range("J113").Value = 701.1429000000001
$ws.Range("K113").Value = 1840.2501
$ws.Range("L113").Value = 2103.4287
$ws.Range("M113").Value = 329.7499
$ws.Range("N113").Value = -6443.4287
$ws.Range("H117").Value = 2409
$ws.Range("I117").Value = 490.5
$ws.Range("J117").Value = 3688
$ws.Range("K117").Value = 1471.5
$ws.Range("L117").Value = 11064
$ws.Range("M117").Value = 1970.5
$ws.Range("N117").Value = -17948
$ws.Range("H129").Value = 1824.8636
$ws.Range("I129").Value = 400
$ws.Range("J129").Value = 1892.7142
$ws.Range("K129").Value = 1200
$ws.Range("L129").Value = 5678.142599999999
$ws.Range("M129").Value = 3800
$ws.Range("N129").Value = -15678.1426
$ws.Range("H131").Value = 802.4343
$ws.Range("I131").Value = 451.8
$ws.Range("J131").Value = 821.0851
$ws.Range("K131").Value = 1355.4
$ws.Range("L131").Value = 2463.2553
$ws.Range("M131").Value = 3684.6
$ws.Range("N131").Value = -12543.2553

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 25000384
$ws.Range("I3").Value = 25000384
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 25000384
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -25000268
$ws.Range("N3").ClearContents()
$ws.Range("H64").Value = 44082
$ws.Range("J64").Value = 44082
$ws.Range("L64").Value = 44082
$ws.Range("N64").Value = -44578
$ws.Range("H67").Value = 44082
$ws.Range("J67").Value = 44082
$ws.Range("L67").Value = 44082
$ws.Range("N67").Value = -45798
$ws.Range("H116").Value = 57900
$ws.Range("J116").Value = 57900
$ws.Range("L116").Value = 57900
$ws.Range("N116").Value = -67078
$ws.Range("H122").Value = 2973.1538
$ws.Range("J122").Value = 2279.8
$ws.Range("L122").Value = 6839.400000000001
$ws.Range("N122").Value = -11739.4
$ws.Range("H132").Value = 2546.1924
$ws.Range("I132").Value = 2329
$ws.Range("J132").Value = 2799.5833
$ws.Range("K132").Value = 6987
$ws.Range("L132").Value = 8398.749899999999
$ws.Range("M132").Value = -4457
$ws.Range("N132").Value = -13458.7499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2046.7273
$ws.Range("I7").Value = 1368.1818
$ws.Range("K7").Value = 1368.1818
$ws.Range("M7").Value = -1256.1818
$ws.Range("H126").Value = 2046.7273
$ws.Range("I126").Value = 1368.1818
$ws.Range("K126").Value = 4104.5454
$ws.Range("M126").Value = -1634.5454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 31710.4
$ws.Range("J56").Value = 34868.777
$ws.Range("L56").Value = 34868.777
$ws.Range("N56").Value = -36296.777
$ws.Range("H122").Value = 1138.1052
$ws.Range("J122").Value = 1385.7142
$ws.Range("L122").Value = 4157.142599999999
$ws.Range("N122").Value = -9057.142599999999
$ws.Range("H126").Value = 1966.0769
$ws.Range("I126").Value = 2340
$ws.Range("J126").Value = 1529.8334
$ws.Range("K126").Value = 7020
$ws.Range("L126").Value = 4589.5002
$ws.Range("M126").Value = -4550
$ws.Range("N126").Value = -9529.5002
